# Update the log write mode: refresh simulated run_time / max_er / iter-N
# columns (rows 2-11) with the latest logged values for gr25_02_simulated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ 3 = 0.2130799293518066; 5 = 28.81181534514144; 6 = 0.0007362173698933311; 7 = 0.0006678772555743052; 8 = 0.0006678772555743052; 9 = 0.0006678772555743052; 10 = 0.0006678772555743052; 11 = 0.0006678772555743052; 12 = 0.0006678772555743052; 13 = 0.0006552036646193739; 14 = 0.0006552036646193739; 15 = 0.0006215379676855998; 16 = 0.0006214602756801723; 17 = 0.0006214602756801723; 18 = 0.0006210546309013405; 19 = 0.0005960337638094946; 20 = 0.0005960337638094946; 21 = 0.0005925026955808834; 22 = 0.0005805270797346169; 23 = 0.0005702430170923889; 24 = 0.0005664384256839151; 25 = 0.0005616338273906712 }
    3 = @{ 3 = 0.2110869884490967; 5 = 28.98144235386644; 7 = 0.0006949423357850833; 8 = 0.0006766368367003506; 9 = 0.0006766368367003506; 10 = 0.0006266364331577898; 11 = 0.0006266364331577898; 12 = 0.0006266364331577898; 13 = 0.0006266364331577898; 14 = 0.0006266364331577898; 15 = 0.0006137435101347228; 16 = 0.0005976461672007483; 17 = 0.0005976461672007483; 18 = 0.0005976461672007483; 19 = 0.0005976461672007483; 20 = 0.0005976461672007483; 21 = 0.0005976461672007483; 22 = 0.0005796493748288302; 23 = 0.0005796493748288302; 24 = 0.0005675286089043292; 25 = 0.0005649403967615289 }
    4 = @{ 3 = 0.2112858295440674; 5 = 28.77917017203981; 7 = 0.0007362173698933311; 8 = 0.0006864051191014704; 9 = 0.0006743744431499749; 10 = 0.0006296718626932314; 11 = 0.0006296718626932314; 12 = 0.0005856721149307188; 13 = 0.0005856721149307188; 14 = 0.0005856721149307188; 15 = 0.0005856721149307188; 16 = 0.0005856721149307188; 17 = 0.0005856721149307188; 18 = 0.0005856721149307188; 19 = 0.0005767475458432153; 20 = 0.0005767475458432153; 21 = 0.0005736662324977824; 22 = 0.0005729544659694248; 23 = 0.0005683373183443077; 24 = 0.000561302373695943; 25 = 0.000560997469240542 }
    5 = @{ 3 = 0.2194716930389404; 5 = 29.15147549329231; 6 = 0.0007362173698933311; 7 = 0.0006915434444275086; 8 = 0.0006915434444275086; 9 = 0.0006915434444275086; 10 = 0.0006915434444275086; 11 = 0.0006915434444275086; 12 = 0.0006743890230094545; 13 = 0.0006515543230889093; 14 = 0.0006483245677867897; 15 = 0.0006331498360020248; 16 = 0.0006331498360020248; 17 = 0.0006169438120222583; 18 = 0.0006082515284678837; 19 = 0.0006075429253827886; 20 = 0.0005787665669873288; 21 = 0.0005787665669873288; 22 = 0.0005787665669873288; 23 = 0.0005752133275916215; 24 = 0.0005706487160214715; 25 = 0.0005682548829101813 }
    6 = @{ 3 = 0.2219798564910889; 5 = 30.03565697046361; 6 = 0.0007357037579825835; 7 = 0.0006756655618497446; 8 = 0.0006646056647404421; 9 = 0.0006646056647404421; 10 = 0.0006646056647404421; 11 = 0.0006646056647404421; 12 = 0.0006646056647404421; 13 = 0.0006646056647404421; 14 = 0.0006541434011217249; 15 = 0.0006541434011217249; 16 = 0.0006508968298264917; 17 = 0.0006218724472918229; 18 = 0.0006218724472918229; 19 = 0.0006218724472918229; 20 = 0.0006218724472918229; 21 = 0.0006089878287543202; 22 = 0.0005989714527548707; 23 = 0.0005880386872435634; 24 = 0.0005880386872435634; 25 = 0.0005854903892877894 }
    7 = @{ 3 = 0.2008368968963623; 5 = 28.96464589875904; 7 = 0.0006898897610080963; 8 = 0.0006898897610080963; 9 = 0.0006445809587068374; 10 = 0.0006445809587068374; 11 = 0.0006445809587068374; 12 = 0.0006445809587068374; 13 = 0.0006445809587068374; 14 = 0.0006445809587068374; 15 = 0.0006312635624574034; 16 = 0.0006152233257688841; 17 = 0.0006152233257688841; 18 = 0.0006152233257688841; 19 = 0.0006009370719071617; 20 = 0.0006004209554164659; 21 = 0.0005912371654597301; 22 = 0.0005768932786495266; 23 = 0.0005746194628148668; 24 = 0.0005667296107469942; 25 = 0.0005646129804826324 }
    8 = @{ 3 = 0.2489378452301025; 5 = 28.37393961933412; 6 = 0.0007362173698933311; 7 = 0.0007362173698933311; 8 = 0.0006688379915143931; 9 = 0.0006688379915143931; 10 = 0.0006659569974699106; 11 = 0.0006659569974699106; 12 = 0.0006648234479981734; 13 = 0.0006468727413231856; 14 = 0.0006323454514170656; 15 = 0.0006323454514170656; 16 = 0.0006175873950330292; 17 = 0.0005910435067072473; 18 = 0.0005910435067072473; 19 = 0.0005831425226375042; 20 = 0.0005791100778199525; 21 = 0.0005634714060556377; 22 = 0.0005634714060556377; 23 = 0.0005634714060556377; 24 = 0.0005539670235583448; 25 = 0.0005530982381936474 }
    9 = @{ 3 = 0.290672779083252; 5 = 29.29229106684943; 6 = 0.0007221195941316113; 7 = 0.0007007931639723821; 8 = 0.0007007931639723821; 9 = 0.0006662518443520462; 10 = 0.0006662518443520462; 11 = 0.0006662518443520462; 12 = 0.0006662518443520462; 13 = 0.0006573641635329469; 14 = 0.0006573641635329469; 15 = 0.0006352158424606226; 16 = 0.000619595567279379; 17 = 0.000619595567279379; 18 = 0.0006030259952790108; 19 = 0.0006030259952790108; 20 = 0.0006019641190567758; 21 = 0.0005951848149244206; 22 = 0.0005853354847700227; 23 = 0.0005841023445775763; 24 = 0.0005790397033194173; 25 = 0.0005709998258645112 }
    10 = @{ 3 = 0.2472438812255859; 5 = 29.34521307516115; 6 = 0.0007362173698933311; 7 = 0.0006910941110176844; 8 = 0.0006910941110176844; 9 = 0.000682868608346538; 10 = 0.000682868608346538; 11 = 0.0006525382891889245; 12 = 0.0006525382891889245; 13 = 0.0006525382891889245; 14 = 0.0006356658704409131; 15 = 0.0006356658704409131; 16 = 0.0006356658704409131; 17 = 0.0006198215652463638; 18 = 0.0006198215652463638; 19 = 0.0006128199539581476; 20 = 0.0006015765827024581; 21 = 0.0005940882775294492; 22 = 0.0005787581277050406; 23 = 0.0005787581277050406; 24 = 0.0005720314439602563; 25 = 0.0005720314439602563 }
    11 = @{ 3 = 0.2373812198638916; 5 = 29.09600839112863; 7 = 0.0007131213600923356; 8 = 0.0006891221997003618; 9 = 0.0006815108183548511; 10 = 0.0006550408121866497; 11 = 0.0006550408121866497; 12 = 0.0006550408121866497; 13 = 0.0006550408121866497; 14 = 0.0006550408121866497; 15 = 0.000650910653675468; 16 = 0.000650910653675468; 17 = 0.0006295866164800543; 18 = 0.0006039191574345961; 19 = 0.0006039191574345961; 20 = 0.0005859808912033545; 21 = 0.0005831964977429984; 22 = 0.0005767384396971443; 23 = 0.0005767384396971443; 24 = 0.0005718892715579599; 25 = 0.0005671736528485113 }
}

foreach ($row in $newValues.Keys) {
    $colMap = $newValues[$row]
    foreach ($col in $colMap.Keys) {
        $ws.Cells.Item($row, $col).Value = $colMap[$col]
    }
}
